$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("таб.8.19")

# Update the absolute path recorded by Excel (cosmetic metadata, best-effort).
# Not exposed via the object model, so this is skipped.

# Move the active selection to L19 on the active sheet.
$ws.Range("L19").Select()

# Add the 2020 column (column M) header and data, matching the existing
# number formatting of the surrounding year columns.
$ws.Range("M4").Value = 2020
$ws.Range("M4").NumberFormat = $ws.Range("L4").NumberFormat

$ws.Range("M5").Value = 34.377950588852634
$ws.Range("M6").Value = 4.8358243107925931
$ws.Range("M7").Value = 5.9543034993102522
$ws.Range("M8").Value = 51.21106605430419
$ws.Range("M9").Value = 27.156801192263725
$ws.Range("M10").Value = 0.94331159862228353
$ws.Range("M11").Value = 7.8509592890793316
$ws.Range("M12").Value = 64.733302669743793
$ws.Range("M13").Value = 97.67954817102779
$ws.Range("M14").Value = 46.725153243037099

$ws.Range("M5:M14").NumberFormat = $ws.Range("L5:L14").NumberFormat

$wb.Save()
